# 10.1.1.xlsx — refresh the year columns from 2007-2019 (16 yrs, cols D:P)
# to 2015-2021 (7 yrs, cols D:J) and drop in the newer data series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up formatting BEFORE the trailing columns disappear ---------------
# Row 5: D5:G5 used a slightly different (legacy) style than H5:J5; make the
# whole row consistent by copying H5's format across to D5:G5.
$ws.Range("H5").Copy()
$ws.Range("D5:G5").PasteSpecial(-4122)

# Row 6: E6:J6 should pick up the style that M6:P6 used to have (D6 keeps
# its original style).
$ws.Range("M6").Copy()
$ws.Range("E6:J6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Drop the now-unused trailing year columns (K:P) ------------------------
$ws.Range("K1:P6").EntireColumn.Delete()

# --- Row 4: year headers ---------------------------------------------------
$years = @(2015, 2016, 2017, 2018, 2019, 2020, 2021)
for ($i = 0; $i -lt $years.Length; $i++) {
    $ws.Cells.Item(4, 4 + $i).Value = $years[$i]
}

# --- Row 5: growth rate among the bottom 40% ------------------------------
$row5 = @(2.2197193775563164, 2.1235271668715399, 2.7818537161298167, 6.7272960584548969, 5.1525830614767187, 4.4774536255935971, 4.6024666695867751)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, 4 + $i).Value = $row5[$i]
}

# --- Row 6: growth rate among the whole population ------------------------
$row6 = @(2.2322863217945752, 2.8603553109638966, 3.113207036164539, 6.2970593463100784, 4.8617746111834492, 2.6715092780025032, 4.3694509108608912)
for ($i = 0; $i -lt $row6.Length; $i++) {
    $ws.Cells.Item(6, 4 + $i).Value = $row6[$i]
}

# --- Column widths: D:J get a uniform custom width ------------------------
$ws.Range("D1:J1").EntireColumn.ColumnWidth = 9.42578125

# --- Selection / cursor position -------------------------------------------
$ws.Range("K16").Select()
